$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-8, columns B through G with new values
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.02744798902934886
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.5159205722746036
$ws.Range("F2").Value = 0.0030670166015625
$ws.Range("G2").Value = 0.9943514917398668

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0.05418734605623229
$ws.Range("D3").Value = 0.1
$ws.Range("E3").Value = 0.0443387513449028
$ws.Range("F3").Value = 0.006764888763427734
$ws.Range("G3").Value = 0.9794267526079573

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 0.05663861903995152
$ws.Range("D4").Value = 0.01
$ws.Range("E4").Value = 0.001107503260070882
$ws.Range("F4").Value = 0.007762908935546875
$ws.Range("G4").Value = 0.975986304295077

$ws.Range("B5").Value = 1885
$ws.Range("C5").Value = 0.05675808151671367
$ws.Range("D5").Value = 0.001
$ws.Range("E5").Value = 0.0009993661110823124
$ws.Range("F5").Value = 1.387660026550293
$ws.Range("G5").Value = 0.6270688390500064

$ws.Range("B6").Value = 10665
$ws.Range("C6").Value = 0.05670706331823926
$ws.Range("D6").Value = 0.0001
$ws.Range("E6").Value = 0.00009959672197821181
$ws.Range("F6").Value = 7.726545810699463
$ws.Range("G6").Value = 0.3707229844828139

$ws.Range("B7").Value = 43163
$ws.Range("C7").Value = 0.05670198295619715
$ws.Range("D7").Value = 0.00001
$ws.Range("E7").Value = 0.000009998218155129855
$ws.Range("F7").Value = 32.39170980453491
$ws.Range("G7").Value = 0.1780661659360412

$ws.Range("B8").Value = 81938
$ws.Range("C8").Value = 0.05670147217535223
$ws.Range("D8").Value = 0.000001
$ws.Range("E8").Value = 0.0000009899626198953622
$ws.Range("F8").Value = 61.82080674171448
$ws.Range("G8").Value = 0.1115551499826239

# Delete row 9 entirely (it no longer exists in the target sheet)
$ws.Range("A9:G9").Delete() | Out-Null
